$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$newValue1 = "U2FsdGVkX19sragcVTigGCm2nUuDT4D+mpoRBq06kxd/HdofgvstuRyBq8MRqt523ciWMWaZJB0fN36DhpLoGbNGBJPfIO2nc3ksb9JkBwS2/sHLN+wzUZzaPFZ3eHCe09myG9/0jSPW4vmFwU0cW/FAy0JenTrvej39oI1H2mMZFy4G8014CTWXMR1d2DXVuKT7vLZfnlJGwcyfpWzN8/wslO2knwakj1pozUB40vARBNaHbkY2eoNG4m3pcdmSdJmCHu++x41rhLkJJBExQMqrUE7jJebu+kLzBIcXlaP05d53NfChweOykwAgXFZXDH8dJ1kcA9ZWvtzrDUcDjRB3zjfj7IcQmh4SOLoI95wZbguQMUgPSH/Xre4IXY8HmIq4g8247SoFUWNt28UZNqmBN9Xy+HCL76TZ5yGtKGATWcS9NaaXPyL9/55WlWdrgyV980ijCExl7HGbVtaluO2ittqgMsc5+0lM1m2Es9bflFrTELggE9O8zLtUQl/Z1uTN3FSvSzXy0b28r9V+PhbY5M54IoggbWwCp6PI9sYktysDCUnIo9l6sFtC86LoxtwVlkoc7Wvskiwo4KRGv7qS/tCUiJbXA/pzj4i/ZYv569zQo9+AQRQsLfYW1SXjYEg16F0eCDzpfuOps9VRria9Aw0pNmRDIlDgPxwvPyimgBnnReTdyFe1W0jJn0zTaTHdJjpOAGzjPaY6UqCVPvj/aPOLmA/topAquKoP6XJ8x3Ug4puYOOAtKrUxc6QA+MiyViLUbndDw874G4ohcjCcfACAl0cZETFVFoJHYq4FpXtO08eeVwjWj9aT6U70tJazSUMxE9fUT73cYiX3iZ+bxXkDrxFD+owxogToT0xwQE/0vut21oU3CBtY0DMxSU9wiqi8NuVu/VtLfL62MB3x5uOg8LfbX9S9g61A5reOB2gleU+aIyjpDdxZGsHzVwI1STkhgxElMpzTaqDQkGbQHMEj0ECCcxvjj1l4TgGHduOvohKagiWV9yp5CUolduU4rYdjFngb+19JWIt1vF4E7RrcYUdzafZLvhaFOOToqfJu/IjrWPnk6xpXwjcZgY2tKzdd/TpV3xemopJlOwozBuhBjCjKuvR/IUAKRN/1eAagMPcFl5br2utSjS7gAFXSkUl1yrHhvkdvndnkVj7OqV7vv2h0yA0MWMGJpPAn6kshhBH5/SvgmT8FIeyL2By/T2mWHpkFCF12F6BdRm8m6/OG0VVU/fZuaaDsvt7KmTurG+RrOdVWH4s7wTXyp1n3HwhKVM7h2ijdumoUYWaH6XGweljaXci5qGX8ijw="

$newValue2 = "U2FsdGVkX1/2fMqQfazvOI1/KlI/nD3tv2y0P97TOx/eG7Z0SfuXnQUk9v1E0z2w5uYXSav66u5dWW0iRSfdniLq/906t4yf/ge+ppabjEWIKk35F8WhUGstKSMOa7vmTUsi/3O3N+Yl+DynVL8cNrr5JXVy1K1tIukx4u/NfEC5/K0kSx5haLCGn20cMSTTiHQwne7ysooJ6W9wltjHGHbyPb4FIK6Qnk5RTowv5KalBzWX+lkomn/qU0i9N9T2t2BbJGqisNgmzQlywtJDZstJa2wLoK57K8aqzH24FeDxa1SjGbuBJSNFUBhKAU0L"

$ws.Range("B2").Value = $newValue1
$ws.Range("B4").Value = $newValue1
$ws.Range("B6").Value = $newValue1

$ws.Range("B3").Value = $newValue2
$ws.Range("B5").Value = $newValue2
$ws.Range("B7").Value = $newValue2
